$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, pushing existing rows 8-16 down to 9-17
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new price-report entry
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C8").Value = "Los Lagos"
$ws.Range("D8").Value = 44495
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 100112013
$ws.Range("G8").Value = "Alcachofa"
$ws.Range("H8").Value = "Madrigal"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 130
$ws.Range("K8").Value = 11000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 11000
$ws.Range("N8").Value = "`$/caja 40 unidades"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 275
$ws.Range("Q8").Value = 40
$ws.Range("R8").Value = "Hortaliza"
